# Actualización automática 2025-10-24 10:30:09
$wb = $excel.ActiveWorkbook

# --- Sheet: VENTAS POR GRUPO ---
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")
$wsGrupo.Range("M31").Value = 3038.01
$wsGrupo.Range("D34").Value = 3232.32
$wsGrupo.Range("M34").Value = 1776.18
$wsGrupo.Range("D41").Value = 1043.04
$wsGrupo.Range("M41").Value = 2391.84

# --- Sheet: VENTA MENSUAL ---
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")
$wsMensual.Range("F31").Value = 5840.61
$wsMensual.Range("F34").Value = 7021.3
$wsMensual.Range("F41").Value = 9587
$wsMensual.Range("F60").Value = 40779.61

# --- Sheet: CUMPLIMIENTO MENSUAL ---
$wsCumpl = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$wsCumpl.Range("D3").Value = 6821.67
$wsCumpl.Range("E3").Value = 13565.8074217135
$wsCumpl.Range("F3").Value = 0.3346009836770998

$wsCumpl.Range("D12").Value = 14001.8
$wsCumpl.Range("E12").Value = 34622.25999999999
$wsCumpl.Range("F12").Value = 0.2879603225234585

$wsCumpl.Range("D14").Value = 44793.55
$wsCumpl.Range("E14").Value = 55104.44284188785
$wsCumpl.Range("F14").Value = 0.4483928928471703
